$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.156.14'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5295'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06323'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07812'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.514'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.679.84'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.884.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8168'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.150.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.600'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.003'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.80%  '
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.212'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("E29").Value = '  +5.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05702'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.276'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.554'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.264'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.596'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.78%  '
$ws.Range("B35").Value = 'MXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.801'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.72%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9509'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5726'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8521'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.807'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.92%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.037.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.797.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈105'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4353'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.857'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05155'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.09%  '
